$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New column header (E1) + data for the week of 21_01_2024
$ws.Range("E1").Value = "21_01_2024"

$ws.Range("E2").Value = 784
$ws.Range("E3").Value = 702
$ws.Range("E4").Value = 1299
$ws.Range("E5").Value = 2814

# Move the active selection like the authored workbook (was A6 -> now E6)
$ws.Range("E6").Select()
